$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1687.5
$ws.Range("I11").Value = 1687.5
$ws.Range("K11").Value = 1687.5
$ws.Range("M11").Value = -1547.5
$ws.Range("H21").Value = 10374.5
$ws.Range("I21").Value = 10374.5
$ws.Range("K21").Value = 10374.5
$ws.Range("M21").Value = -9906.5
$ws.Range("H23").Value = 10374.5
$ws.Range("I23").Value = 10374.5
$ws.Range("K23").Value = 10374.5
$ws.Range("M23").Value = -10140.5
$ws.Range("H40").Value = 5487.5
$ws.Range("I40").Value = 5557.143
$ws.Range("K40").Value = 5557.143
$ws.Range("M40").Value = -5382.143
$ws.Range("H70").Value = 6315.9565
$ws.Range("I70").Value = 5944.25
$ws.Range("J70").Value = 6394.2104
$ws.Range("K70").Value = 17832.75
$ws.Range("L70").Value = 19182.6312
$ws.Range("M70").Value = -17562.75
$ws.Range("N70").Value = -19722.6312
$ws.Range("H73").Value = 6315.9565
$ws.Range("I73").Value = 5944.25
$ws.Range("J73").Value = 6394.2104
$ws.Range("K73").Value = 17832.75
$ws.Range("L73").Value = 19182.6312
$ws.Range("M73").Value = -16896.75
$ws.Range("N73").Value = -21054.6312
$ws.Range("H103").Value = 609.3
$ws.Range("I103").Value = 671
$ws.Range("J103").Value = 362.5
$ws.Range("K103").Value = 2013
$ws.Range("L103").Value = 1087.5
$ws.Range("M103").Value = -1427
$ws.Range("N103").Value = -2259.5
$ws.Range("H111").Value = 10104488
$ws.Range("I111").Value = 12348771
$ws.Range("J111").Value = 5216
$ws.Range("K111").Value = 37046313
$ws.Range("L111").Value = 15648
$ws.Range("M111").Value = -37043246
$ws.Range("N111").Value = -21782
$ws.Range("H125").Value = 9012035
$ws.Range("I125").Value = 1514
$ws.Range("J125").Value = 12824179
$ws.Range("K125").Value = 13626
$ws.Range("L125").Value = 115417611
$ws.Range("M125").Value = -11166
$ws.Range("N125").Value = -115422531
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3751.79
$ws.Range("I32").Value = 2065.2705
$ws.Range("J32").Value = 13308.733
$ws.Range("K32").Value = 2065.2705
$ws.Range("L32").Value = 13308.733
$ws.Range("M32").Value = -1778.2705
$ws.Range("N32").Value = -13882.733
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H122").Value = 1306080.5
$ws.Range("I122").Value = 3167.6667
$ws.Range("K122").Value = 9503.000100000001
$ws.Range("M122").Value = -7053.000100000001
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8406738
$ws.Range("I99").Value = 17859636
$ws.Range("J99").Value = 4163.3335
$ws.Range("K99").Value = 17859636
$ws.Range("L99").Value = 4163.3335
$ws.Range("M99").Value = -17858138
$ws.Range("N99").Value = -7159.3335
$ws.Range("H103").Value = 27884.666
$ws.Range("J103").Value = 27884.666
$ws.Range("L103").Value = 27884.666
$ws.Range("N103").Value = -30228.666
$ws.Range("H107").Value = 10206139
$ws.Range("I107").Value = 11906446
$ws.Range("J107").Value = 4300
$ws.Range("K107").Value = 11906446
$ws.Range("L107").Value = 4300
$ws.Range("M107").Value = -11904526
$ws.Range("N107").Value = -8140
$ws.Range("H134").Value = 3451.9736
$ws.Range("I134").Value = 1214.1
$ws.Range("K134").Value = 3642.3
$ws.Range("M134").Value = -1107.3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2317.1667
$ws.Range("I16").Value = 1980.6
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1980.6
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -1693.6
$ws.Range("N16").Value = -4574
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H113").Value = 2317.1667
$ws.Range("I113").Value = 1980.6
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1980.6
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 189.4000000000001
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 1699.0975
$ws.Range("I122").Value = 1620.2
$ws.Range("J122").Value = 1914.2727
$ws.Range("K122").Value = 4860.6
$ws.Range("L122").Value = 5742.8181
$ws.Range("M122").Value = -2410.6
$ws.Range("N122").Value = -10642.8181
$ws.Range("H141").Value = 171429.69
$ws.Range("J141").Value = 171429.69
$ws.Range("L141").Value = 171429.69
$ws.Range("N141").Value = -181789.69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 59613.47
$ws.Range("I5").Value = 764.8461
$ws.Range("J5").Value = 250871.5
$ws.Range("K5").Value = 2294.5383
$ws.Range("L5").Value = 752614.5
$ws.Range("M5").Value = -2182.5383
$ws.Range("N5").Value = -752838.5
$ws.Range("H38").Value = 62.944443
$ws.Range("I38").Value = 25.818182
$ws.Range("J38").Value = 121.28571
$ws.Range("K38").Value = 77.45454599999999
$ws.Range("L38").Value = 363.85713
$ws.Range("M38").Value = 269.545454
$ws.Range("N38").Value = -1057.85713
$ws.Range("H68").Value = 1489.875
$ws.Range("J68").Value = 2863
$ws.Range("L68").Value = 8589
$ws.Range("N68").Value = -10211
$ws.Range("H71").Value = 1489.875
$ws.Range("J71").Value = 2863
$ws.Range("L71").Value = 25767
$ws.Range("N71").Value = -33879
$ws.Range("H92").Value = 670
$ws.Range("J92").Value = 2710
$ws.Range("L92").Value = 8130
$ws.Range("N92").Value = -10626
$ws.Range("H94").Value = 9922.182000000001
$ws.Range("I94").Value = 9624.5
$ws.Range("J94").Value = 9988.333000000001
$ws.Range("K94").Value = 28873.5
$ws.Range("L94").Value = 29964.999
$ws.Range("M94").Value = -28197.5
$ws.Range("N94").Value = -31316.999
$ws.Range("H107").Value = 296.25
$ws.Range("I107").Value = 267.14285
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 801.4285500000001
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1118.57145
$ws.Range("N107").Value = -5340
$ws.Range("H131").Value = 14371662
$ws.Range("I131").Value = 10419631
$ws.Range("J131").Value = 15877198
$ws.Range("K131").Value = 31258893
$ws.Range("L131").Value = 47631594
$ws.Range("M131").Value = -31253853
$ws.Range("N131").Value = -47641674
$ws.Range("H132").Value = 2697
$ws.Range("I132").Value = 1498.75
$ws.Range("J132").Value = 3039.3572
$ws.Range("K132").Value = 13488.75
$ws.Range("L132").Value = 27354.2148
$ws.Range("M132").Value = -10958.75
$ws.Range("N132").Value = -32414.2148
$ws.Range("H135").Value = 59613.47
$ws.Range("I135").Value = 764.8461
$ws.Range("J135").Value = 250871.5
$ws.Range("K135").Value = 6883.6149
$ws.Range("L135").Value = 2257843.5
$ws.Range("M135").Value = -4348.6149
$ws.Range("N135").Value = -2262913.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 31388
$ws.Range("J32").Value = 31388
$ws.Range("L32").Value = 31388
$ws.Range("N32").Value = -31980
$ws.Range("H97").Value = 2382148.2
$ws.Range("I97").Value = 3402226.2
$ws.Range("K97").Value = 3402226.2
$ws.Range("M97").Value = -3401730.2
$ws.Range("H102").Value = 4690478
$ws.Range("I102").Value = 15875339
$ws.Range("K102").Value = 15875339
$ws.Range("M102").Value = -15873717
$ws.Range("H107").Value = 861.4
$ws.Range("I107").Value = 1101.3334
$ws.Range("K107").Value = 1101.3334
$ws.Range("M107").Value = 818.6666
$ws.Range("H109").Value = 49362.25
$ws.Range("J109").Value = 49362.25
$ws.Range("L109").Value = 49362.25
$ws.Range("N109").Value = -51442.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 39999
$ws.Range("I41").Value = 39999
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 39999
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -39561
$ws.Range("H93").Value = 9531260
$ws.Range("I93").Value = 13335153
$ws.Range("J93").Value = 21528.7
$ws.Range("K93").Value = 13335153
$ws.Range("L93").Value = 21528.7
$ws.Range("M93").Value = -13333905
$ws.Range("N93").Value = -24024.7
$ws.Range("H100").Value = 3291.7083
$ws.Range("I100").Value = 2763.2632
$ws.Range("K100").Value = 2763.2632
$ws.Range("M100").Value = -2222.2632
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 25000
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -29992
$ws.Range("H94").Value = 21247
$ws.Range("J94").Value = 24996
$ws.Range("L94").Value = 24996
$ws.Range("N94").Value = -26798
$ws.Range("H100").Value = 1487.5454
$ws.Range("I100").Value = 2116.3333
$ws.Range("K100").Value = 4232.6666
$ws.Range("M100").Value = -3691.6666
$ws.Range("H107").Value = 55556984
$ws.Range("I107").Value = 111111624
$ws.Range("J107").Value = 2348.889
$ws.Range("K107").Value = 333334872
$ws.Range("L107").Value = 7046.667
$ws.Range("M107").Value = -333332952
$ws.Range("N107").Value = -10886.667
$ws.Range("H126").Value = 3629.8
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H129").Value = 58427.2
$ws.Range("J129").Value = 75281.60000000001
$ws.Range("L129").Value = 75281.60000000001
$ws.Range("N129").Value = -85281.60000000001
